$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.7020156925206528
$ws.Range("J2").Value = 0.7020156925206527
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 0.279882508464
$ws.Range("R2").Value = 2.518942576176
$ws.Range("S2").Value = 0.007906448660253503
$ws.Range("T2").Value = 0.007906448660253503

$ws.Range("I3").Value = 0.7020156925206528
$ws.Range("J3").Value = 0.7020156925206527
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("S3").Value = 0.6445245137905832
$ws.Range("T3").Value = 0.6445245137905832

$ws.Range("I4").Value = 0.7020156925206528
$ws.Range("J4").Value = 0.7020156925206527
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 1.755263232558
$ws.Range("R4").Value = 15.797369093022
$ws.Range("S4").Value = 0.04958473006981599
$ws.Range("T4").Value = 0.04958473006981599

$ws.Range("G5").Value = 1.074622
$ws.Range("H5").Value = 3.223866
$ws.Range("I5").Value = 0.2979843074793473
$ws.Range("J5").Value = 0.2979843074793473
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 0.118801611344
$ws.Range("R5").Value = 1.069214502096
$ws.Range("S5").Value = 0.00335604695699497
$ws.Range("T5").Value = 0.00335604695699497

$ws.Range("G6").Value = 1.074622
$ws.Range("H6").Value = 3.223866
$ws.Range("I6").Value = 0.2979843074793473
$ws.Range("J6").Value = 0.2979843074793473
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("Q6").Value = 9.684569403954667
$ws.Range("R6").Value = 87.161124635592
$ws.Range("S6").Value = 0.2735810508818501
$ws.Range("T6").Value = 0.2735810508818502

$ws.Range("G7").Value = 1.074622
$ws.Range("H7").Value = 3.223866
$ws.Range("I7").Value = 0.2979843074793473
$ws.Range("J7").Value = 0.2979843074793473
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 0.7450558504179999
$ws.Range("R7").Value = 6.705502653761999
$ws.Range("S7").Value = 0.02104720964050217
$ws.Range("T7").Value = 0.02104720964050217
